$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix GCRMN ecoregion assignments for ROPME / PERSGA (rows 64-69).
# Column C (MEOW Ecoregion) values stay the same; the GCRMN Region (A)
# and Subregion (B) labels were swapped between the PERSGA and ROPME
# blocks.

$ws.Cells.Item(64, 1).Value = "ROPME"
$ws.Cells.Item(64, 2).Value = "ROPME.1"

$ws.Cells.Item(65, 1).Value = "ROPME"
$ws.Cells.Item(65, 2).Value = "ROPME.2"

$ws.Cells.Item(66, 1).Value = "ROPME"
$ws.Cells.Item(66, 2).Value = "ROPME.3"

$ws.Cells.Item(67, 1).Value = "PERSGA"
$ws.Cells.Item(67, 2).Value = "PERSGA.1"

$ws.Cells.Item(68, 1).Value = "PERSGA"
$ws.Cells.Item(68, 2).Value = "PERSGA.2"

$ws.Cells.Item(69, 1).Value = "PERSGA"
$ws.Cells.Item(69, 2).Value = "PERSGA.3"

# Update the view state: scroll the frozen pane and land the selection on
# G53 (matches the workbook's sheetView/selection in the target revision).
$ws.Range("A23").Select()
$ws.Range("G53").Select()
